# Fixing map precision + Scenes in positive octant
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content / text fixes
$ws.Range("A35").Value = "Crytek Sponza (pre ASTC)"
$ws.Range("E2").Value = "256x256 is worth for Sponza"
$ws.Range("D8").Value = "DEFERRED NAIVE FPS"
$ws.Range("E8").Value = "DEFERRED INTERPOLATED FPS"

# Widen column E to fit the longer label
$ws.Columns.Item(5).ColumnWidth = 28.6

# Update the saved view state (scroll position + active selection)
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E37").Select()
